$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N (14), shifting N:P -> O:Q
$ws.Columns.Item(14).Insert()

# Match the target column width for the newly inserted column (raw XML width 11)
$ws.Columns.Item(14).ColumnWidth = 10.17

# Make "Repayment schedule" the active sheet / tab, with the new selection
$ws.Activate()
$ws.Range("R8").Select()
